$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New measurement rows (44-69) collected for Spratelloides_gracilis on a
# second collection date (5_16_2023). Columns:
#   row, Mass_g(D), SL_mm(E), TL_mm(F), Lot_USNM_ID(J)
# A (Species) = "Spratelloides_gracilis", B (Species_nu) = 9,
# C (Individual_#) = row - 1, I (Date) = "5_16_2023"
# ---------------------------------------------------------------------------
$newRows = @(
  @(44, 0.23300000000000001, 32.4, 35, 120834),
  @(45, 0.39, 35.5, 37.700000000000003, 120834),
  @(46, 0.36, 35.200000000000003, 40, 120834),
  @(47, 0.373, 36.5, 39.4, 120834),
  @(48, 0.15, 28.8, 31.1, 120834),
  @(49, 0.39500000000000002, 34.799999999999997, 37.5, 120834),
  @(50, 0.25, 34.4, 37.5, 120834),
  @(51, 0.29399999999999998, 32.799999999999997, 35.299999999999997, 120834),
  @(52, 0.26800000000000002, 33.700000000000003, 36.799999999999997, 120834),
  @(53, 0.246, 31.9, 34.1, 120834),
  @(54, 0.34200000000000003, 35.4, 38.299999999999997, 120834),
  @(55, 0.38500000000000001, 36.9, 40.5, 120834),
  @(56, 0.32500000000000001, 33.5, 36.700000000000003, 120834),
  @(57, 0.23100000000000001, 33.200000000000003, 36.5, 120834),
  @(58, 0.36599999999999999, 36, 40.1, 120834),
  @(59, 0.28199999999999997, 34.1, 37, 120834),
  @(60, 0.28499999999999998, 33.5, 36.700000000000003, 120834),
  @(61, 0.255, 34.299999999999997, 36.4, 120834),
  @(62, 0.42899999999999999, 41.1, 44.2, 120834),
  @(63, 0.374, 37, 41.2, 120834),
  @(64, 0.29499999999999998, 34.4, 37.5, 120834),
  @(65, 0.08, 23.1, 25, 120834),
  @(66, 0.151, 27.9, 30, 120834),
  @(67, 0.3, 35.799999999999997, 38.200000000000003, 120834),
  @(68, 0.33900000000000002, 37.9, 41.6, 120834),
  @(69, 0.38300000000000001, 37, 41, 120834)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $mass = $r[1]
    $sl = $r[2]
    $tl = $r[3]
    $lot = $r[4]

    $ws.Range("A$row").Value2 = "Spratelloides_gracilis"
    $ws.Range("B$row").Value2 = 9
    $ws.Range("C$row").Value2 = ($row - 1)
    $ws.Range("D$row").Value2 = $mass
    $ws.Range("E$row").Value2 = $sl
    $ws.Range("F$row").Value2 = $tl
    $ws.Range("I$row").Value2 = "5_16_2023"
    $ws.Range("J$row").Value2 = $lot
}

# ---------------------------------------------------------------------------
# Formulas: extend the existing shared formula (SL_cm/TL_cm = SL_mm/TL_mm /10)
# down through the new rows. The original block G34:H43 grows to G34:H49,
# and a further block G50:H71 is created (matching the larger selection the
# author made while filling down before only populating data to row 69).
# ---------------------------------------------------------------------------
$ws.Range("G34:H49").Formula = "=E34/10"
$ws.Range("G50:H71").Formula = "=E50/10"
# Rows 70:71 were never populated with data, so remove those trailing,
# otherwise-empty rows again (keeps the shared-formula ref spanning to 71
# without leaving stray empty rows behind).
$ws.Range("G70:H71").ClearContents()

# ---------------------------------------------------------------------------
# Formatting: reuse the workbook's existing cell styles (index 2 = 10pt
# Calibri, index 3 = 10pt Arial) for the appropriate columns/rows, exactly
# as the alternating pattern already used in rows 2-43 continues.
# ---------------------------------------------------------------------------
foreach ($r in $newRows) {
    $row = $r[0]
    $isEven = ($row % 2) -eq 0

    # B:H always use style 2
    $ws.Range("B2").Copy() | Out-Null
    $ws.Range("B$row`:H$row").PasteSpecial(-4122) | Out-Null

    if ($isEven) {
        # A uses style 2, I uses style 3, J uses style 2
        $ws.Range("B2").Copy() | Out-Null
        $ws.Range("A$row").PasteSpecial(-4122) | Out-Null

        $ws.Range("I2").Copy() | Out-Null
        $ws.Range("I$row").PasteSpecial(-4122) | Out-Null

        $ws.Range("J2").Copy() | Out-Null
        $ws.Range("J$row").PasteSpecial(-4122) | Out-Null
    } else {
        # A uses style 3; I and J are left unformatted (no explicit style)
        $ws.Range("A3").Copy() | Out-Null
        $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    }
}

# ---------------------------------------------------------------------------
# Update the active selection to reflect where the author left off editing.
# ---------------------------------------------------------------------------
$ws.Range("A2:J71").Select() | Out-Null
